# Apply the commit's changes to the PlayerPerformance workbook:
#  1. Insert a new "Player Info" worksheet in front of the existing sheets,
#     with player ID/name/batting-hand/bowling-style data.
#  2. Rename the MATCH_CARD_LINK column to MATCH_CODE on both the
#     "ODI Batting" and "ODI Bowling" sheets, and replace the full
#     scorecard URL value with just the numeric match code.

$wb = $excel.ActiveWorkbook

# --- 1. Insert new "Player Info" worksheet at the front of the workbook ---
$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet)
$playerInfo.Name = "Player Info"

# Header row: bold, centered, top-aligned, thin border (matches the style
# already used for header rows on the other sheets in this workbook).
$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop
$headerRange.Borders.LineStyle = 1

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Data row. The leading "'" forces the numeric-looking ID to be stored as
# text (matching the original inlineStr type), and re-applying the
# "Normal" style afterwards clears the quote-prefix formatting so the
# cell keeps the plain, unstyled look of the other data rows.
$playerInfo.Range("A2").Value = "'5984"
$playerInfo.Range("A2").Style = "Normal"
$playerInfo.Range("B2").Value = "Lizaad Buyron Williams"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium Fast"

# --- 2. "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE ---
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("D1").Value = "MATCH_CODE"
$odiBatting.Range("D2").Value = "'4478"
$odiBatting.Range("D2").Style = "Normal"

# --- 3. "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE ---
$odiBowling = $wb.Worksheets.Item("ODI Bowling")
$odiBowling.Range("B1").Value = "MATCH_CODE"
$odiBowling.Range("B2").Value = "'4478"
$odiBowling.Range("B2").Style = "Normal"
